# repull data, push all data, mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "dSF" column (F) values to reflect the repulled data
$ws.Range("F2").Value = -1
$ws.Range("F4").Value = 8
$ws.Range("F5").Value = -2
$ws.Range("F9").Value = -1
$ws.Range("F10").Value = 0
$ws.Range("F17").Value = 4
